# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values on row 4 of the
# zh-cn and de-de sheets to reflect the new report generation times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-14 09:09:35"
$wsZhCn.Range("H4").Value = "2016-03-14 09:10:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-14 09:09:44"
$wsDeDe.Range("H4").Value = "2016-03-14 09:10:43"
